$d = $word.ActiveDocument

function Set-ParagraphText($oldText, $newText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($oldText + "`r")) {
            $p.Range.Text = $newText
            return $true
        }
    }
    return $false
}

function Remove-ParagraphByText($matchText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($matchText + "`r")) {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

# Heading line
Set-ParagraphText 'New in this update (Railway frontend dependency fix)' 'New in this update (Render deployment setup)'

# "New in this update" bullet block rewrite (5 lines 1:1)
Set-ParagraphText '- Fixed frontend Docker build failure due to peer dependency conflict:' '- Added Render Blueprint file: `render.yaml`.'
Set-ParagraphText '  - `react-day-picker@8.10.1` expects `date-fns ^2/^3`, while project uses `date-fns@4`.' '  - Provisions PostgreSQL, backend service, and frontend static service.'
Set-ParagraphText '- Updated `frontend/Dockerfile` install command to:' '- Added Render deploy runbook: `DEPLOY_RENDER.md`.'
Set-ParagraphText '  - `npm install --legacy-peer-deps`' '- Render build config includes frontend install fallback:'
Set-ParagraphText '- This resolves Railway build error `ERESOLVE unable to resolve dependency tree`.' '  - `npm install --legacy-peer-deps && npm run build`.'

# Remove the whole "Proper long-term dependency alignment..." paragraph
Remove-ParagraphByText '- Proper long-term dependency alignment still pending (recommended: migrate `date-fns` to a version compatible with all peers or upgrade dependent packages).'

# Git state updates
Set-ParagraphText '- Last pushed commit: 9cfabad' '- Last pushed commit: 09217f9'
Set-ParagraphText '- Current frontend dependency build fix is local and not pushed yet.' '- Current Render deployment setup is local and not pushed yet.'
